# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-job Leve
# profit sheets. Values come from an external price-fetch pass, so
# they're written as plain literals (no formulas involved).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4024.25
$ws.Range("I9").Value = 199.66667
$ws.Range("K9").Value = 199.66667
$ws.Range("M9").Value = -30.66667000000001

$ws.Range("H33").Value = 40423.12
$ws.Range("I33").Value = 52857.05
$ws.Range("K33").Value = 52857.05
$ws.Range("M33").Value = -52628.05

$ws.Range("H74").Value = 6155.1333
$ws.Range("I74").Value = 5952
$ws.Range("K74").Value = 5952
$ws.Range("M74").Value = -5016

$ws.Range("H77").Value = 6155.1333
$ws.Range("I77").Value = 5952
$ws.Range("K77").Value = 29760
$ws.Range("M77").Value = -25080

$ws.Range("H92").Value = 779.2857
$ws.Range("J92").Value = 671.5
$ws.Range("L92").Value = 671.5
$ws.Range("N92").Value = -3167.5

$ws.Range("H94").Value = 999
$ws.Range("I94").Value = 999
$ws.Range("K94").Value = 999
$ws.Range("M94").Value = -548

$ws.Range("H100").Value = 5052.5625
$ws.Range("J100").Value = 4374.75
$ws.Range("L100").Value = 4374.75
$ws.Range("N100").Value = -5456.75

$ws.Range("H107").Value = 1528.3
$ws.Range("I107").Value = 1587.1111
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 1587.1111
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 332.8888999999999
$ws.Range("N107").Value = -4839

$ws.Range("H111").Value = 2821.3333
$ws.Range("I111").Value = 2464
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 7392
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -4325
$ws.Range("N111").Value = -15134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3408.5
$ws.Range("J2").Value = 3666.3333
$ws.Range("L2").Value = 3666.3333
$ws.Range("N2").Value = -3892.3333

$ws.Range("H4").Value = 175.66667
$ws.Range("I4").Value = 191
$ws.Range("K4").Value = 191
$ws.Range("M4").Value = -75

$ws.Range("H32").Value = 3476.5527
$ws.Range("I32").Value = 3476.5527
$ws.Range("K32").Value = 3476.5527
$ws.Range("M32").Value = -3189.5527

$ws.Range("H110").Value = 1598.625
$ws.Range("I110").Value = 1535.2
$ws.Range("K110").Value = 1535.2
$ws.Range("M110").Value = 509.8

$ws.Range("H116").Value = 3408.5
$ws.Range("J116").Value = 3666.3333
$ws.Range("L116").Value = 3666.3333
$ws.Range("N116").Value = -8254.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3408.5
$ws.Range("J3").Value = 3666.3333
$ws.Range("L3").Value = 3666.3333
$ws.Range("N3").Value = -3894.3333

$ws.Range("H36").Value = 7824.75
$ws.Range("I36").Value = 2099.6667
$ws.Range("J36").Value = 25000
$ws.Range("K36").Value = 2099.6667
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = -1565.6667
$ws.Range("N36").Value = -26068

$ws.Range("H107").Value = 2744
$ws.Range("J107").Value = 1899.5
$ws.Range("L107").Value = 1899.5
$ws.Range("N107").Value = -5739.5

$ws.Range("H135").Value = 19999
$ws.Range("J135").Value = 19999
$ws.Range("L135").Value = 19999
$ws.Range("N135").Value = -30139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 11861.556
$ws.Range("I33").Value = 7881.25
$ws.Range("J33").Value = 15045.8
$ws.Range("K33").Value = 7881.25
$ws.Range("L33").Value = 15045.8
$ws.Range("M33").Value = -7502.25
$ws.Range("N33").Value = -15803.8

$ws.Range("H50").Value = 18250
$ws.Range("J50").Value = 18250
$ws.Range("L50").Value = 18250
$ws.Range("N50").Value = -19500

$ws.Range("H59").Value = 23119.334
$ws.Range("J59").Value = 25127
$ws.Range("L59").Value = 25127
$ws.Range("N59").Value = -27417

$ws.Range("H60").Value = 11997.4
$ws.Range("I60").Value = 9327.666999999999
$ws.Range("J60").Value = 16002
$ws.Range("K60").Value = 9327.666999999999
$ws.Range("L60").Value = 16002
$ws.Range("M60").Value = -8816.666999999999
$ws.Range("N60").Value = -17024

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 168092.06
$ws.Range("I4").Value = 221.58333
$ws.Range("K4").Value = 664.74999
$ws.Range("M4").Value = -552.74999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1299.375
$ws.Range("I22").Value = 2065.3333
$ws.Range("J22").Value = 839.8
$ws.Range("K22").Value = 2065.3333
$ws.Range("L22").Value = 839.8
$ws.Range("M22").Value = -1770.3333
$ws.Range("N22").Value = -1429.8

$ws.Range("H27").Value = 1299.375
$ws.Range("I27").Value = 2065.3333
$ws.Range("J27").Value = 839.8
$ws.Range("K27").Value = 2065.3333
$ws.Range("L27").Value = 839.8
$ws.Range("M27").Value = -1958.3333
$ws.Range("N27").Value = -1053.8

$ws.Range("H68").Value = 3684
$ws.Range("I68").Value = 3684
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3684
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2935
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 3684
$ws.Range("I71").Value = 3684
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 18420
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -14676
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 2873.4707
$ws.Range("I82").Value = 3349.8462
$ws.Range("K82").Value = 3349.8462
$ws.Range("M82").Value = -2988.8462

$ws.Range("H85").Value = 2873.4707
$ws.Range("I85").Value = 3349.8462
$ws.Range("K85").Value = 3349.8462
$ws.Range("M85").Value = -2101.8462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4571.4287
$ws.Range("I62").Value = 4863.636
$ws.Range("K62").Value = 4863.636
$ws.Range("M62").Value = -4239.636

$ws.Range("H65").Value = 4571.4287
$ws.Range("I65").Value = 4863.636
$ws.Range("K65").Value = 24318.18
$ws.Range("M65").Value = -21198.18

$ws.Range("H113").Value = 935.1
$ws.Range("I113").Value = 935.1
$ws.Range("K113").Value = 2805.3
$ws.Range("M113").Value = -635.3000000000002

$ws.Range("H126").Value = 3099.353
$ws.Range("I126").Value = 2799.375
$ws.Range("J126").Value = 7899
$ws.Range("K126").Value = 8398.125
$ws.Range("L126").Value = 23697
$ws.Range("M126").Value = -5928.125
$ws.Range("N126").Value = -28637
